$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.500406265258789
$ws.Range("B1").Value = 1.527007341384888
$ws.Range("C1").Value = 1.713122367858887
$ws.Range("D1").Value = 2.634259223937988
$ws.Range("E1").Value = -1
